$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldVal = "dnasr281@gmail.com, System"
$newVal = "System, dnasr281@gmail.com"

$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count

$changed = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2
    if ($v -eq $oldVal) {
        $cell.Value2 = $newVal
        $changed = $changed + 1
    }
}

Write-Host "Cells changed: " $changed
